# Automatische test-sync: 2025-06-29 15:24:50
#
# Adds a new "Testmail #20" row to the Logs sheet and refreshes the
# Dashboard category-count summary (which is sorted by count, so the
# "Overig" / "Retour / Terugbetaling" rows swap places once "Overig"
# overtakes it).

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Logs: append row 35 -------------------------------------------------
$newRow = 35

$logs.Cells.Item($newRow, 1).Value = "Kun je deze taak op je nemen?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #20: Kun je deze taak op je nemen?"
$logs.Cells.Item($newRow, 4).Value = "Overig"
$logs.Cells.Item($newRow, 5).Value = "Beste [Naam],`nBedankt voor je bericht. Kun je wat meer details geven over welke taak je precies bedoelt? Dan kan ik je zo goed mogelijk helpen.`nMet vriendelijke groet,`n[Jouw naam]  `nE-mailassistent bij [Bedrijfsnaam]"
$logs.Cells.Item($newRow, 6).Value = "2025-06-29 15:24:39"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Ja"
$logs.Cells.Item($newRow, 9).Value = "Nee"

# --- Dashboard: "Overig" now has 4 hits, overtaking "Retour / Terugbetaling" (3) ---
$dash.Cells.Item(5, 1).Value = "Overig"
$dash.Cells.Item(5, 2).Value = 4
$dash.Cells.Item(6, 1).Value = "Retour / Terugbetaling"
$dash.Cells.Item(6, 2).Value = 3

# --- Logs: extend the conditional-formatting ranges to cover the new row ---
foreach ($col in @("D", "G", "H", "I")) {
    $oldRange = $logs.Range("$col`2:$col`34")
    $newRange = $logs.Range("$col`2:$col`35")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}
